$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("externalShortName", $false, $false, $false, $false, $false, `
              $true, 1, $false, "external_short_name", 2)
